$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, as scraped from the source diff.
$updates = [ordered]@{
    "D2" = "60.767.30"
    "E2" = "  -1.33%  "
    "D3" = "3.383.78"
    "E3" = "  -1.95%  "
    "D4" = "1.00"
    "E4" = "  +0.01%  "
    "D5" = "569.13"
    "E5" = "  -1.71%  "
    "D6" = "140.91"
    "E6" = "  -2.66%  "
    "E7" = "  +0.02%  "
    "D8" = "3.383.84"
    "E8" = "  -1.99%  "
    "E9" = "  -0.48%  "
    "D10" = "7.51"
    "E10" = "  -1.30%  "
    "E11" = "  -1.56%  "
    "D12" = "0.394"
    "E12" = "  +1.68%  "
    "D13" = "3.960.77"
    "E13" = "  -1.97%  "
    "D14" = "28.46"
    "E15" = "  +2.22%  "
    "E16" = "  -1.52%  "
    "D17" = "3.380.62"
    "E17" = "  -1.73%  "
    "D18" = "60.849.63"
    "E18" = "  -1.41%  "
    "D19" = "6.20"
    "E19" = "  -0.87%  "
    "E20" = "  -1.99%  "
    "D22" = "382.88"
    "E22" = "  -1.68%  "
    "D23" = "0.558"
    "E23" = "  -1.00%  "
    "D24" = "73.67"
    "E24" = "  +0.41%  "
    "E25" = "  +0.42%  "
    "D26" = "0.0000116"
    "E26" = "  -5.35%  "
    "D27" = "3.518.84"
    "E27" = "  -1.94%  "
    "E28" = "  -0.32%  "
    "D29" = "0.999"
    "E29" = "  -0.28%  "
    "D30" = "7.38"
    "D31" = "7.98"
    "E31" = "  -1.98%  "
    "E32" = "  -1.81%  "
    "E33" = "  -2.81%  "
    "E35" = "  -1.83%  "
    "D36" = "6.95"
    "D37" = "166.29"
    "E37" = "  -0.39%  "
    "D38" = "3.413.90"
    "E38" = "  -1.89%  "
    "D39" = "4.97"
    "E39" = "  -3.01%  "
    "E40" = "  -4.55%  "
    "D41" = "27.92"
    "E41" = "  -0.57%  "
    "D42" = "0.0771"
    "E42" = "  -1.00%  "
    "E43" = "  +0.02%  "
    "D44" = "0.777"
    "E44" = "  -2.80%  "
    "D45" = "41.90"
    "E45" = "  -0.97%  "
    "D46" = "4.40"
    "E46" = "  -1.67%  "
    "E47" = "  -3.83%  "
    "D48" = "1.13"
    "D49" = "2.491.90"
    "E49" = "  -3.57%  "
    "D50" = "23.58"
    "E50" = "  +2.64%  "
    "D51" = "6.80"
    "E51" = "  -1.60%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text formatting so purely-numeric-looking strings (e.g. "1.00", "7.98")
    # are preserved exactly as-is instead of being normalised to a number by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
